$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Updated product backlog priorities for the first few user stories
$ws.Range("A2").Value = 0.5
$ws.Range("A3").Value = 0.5
$ws.Range("A4").Value = 1

# Reflect the editor's current view (zoom level + selected cell) when saved
$excel.ActiveWindow.Zoom = 129
$ws.Range("A4").Select()
